$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -0.318
$ws.Range("D3").Value = -0.318
$ws.Range("G2").Value = -0.7697368421052631
$ws.Range("G3").Value = -0.7697368421052631
$ws.Range("H2").Value = -0.7697368421052631
$ws.Range("H3").Value = -0.7697368421052631
$ws.Range("I2").Value = -0.5805612910675676
$ws.Range("I3").Value = -0.5805612910675676
$ws.Range("J2").Value = -0.5805612910675676
$ws.Range("J3").Value = -0.5805612910675676
$ws.Range("K2").Value = -2.66
$ws.Range("K3").Value = -2.66
$ws.Range("L2").Value = -1.75
$ws.Range("L3").Value = -1.75
$ws.Range("U2").Value = 0.749
$ws.Range("U3").Value = 0.749
$ws.Range("V2").Value = 0.3242424242424242
$ws.Range("V3").Value = 0.3242424242424242
$ws.Range("W2").Value = -0.6734177215189874
$ws.Range("W3").Value = -0.6734177215189874
$ws.Range("X2").Value = 0.07460332883262136
$ws.Range("X3").Value = 0.07460332883262136
$ws.Range("Y2").Value = -0.7480210503516087
$ws.Range("Y3").Value = -0.7480210503516087
$ws.Range("Z2").Value = 0.4401630464321856
$ws.Range("Z3").Value = 0.4401630464321856
$ws.Range("AA2").Value = -0.2555416265169034
$ws.Range("AA3").Value = -0.2555416265169034
$ws.Range("AB2").Value = 0.07077732592338744
$ws.Range("AB3").Value = 0.07077732592338744
$ws.Range("AC2").Value = -0.3263189524402908
$ws.Range("AC3").Value = -0.3263189524402908
$ws.Range("AD2").Value = 0.286
$ws.Range("AD3").Value = 0.286
$ws.Range("AE2").Value = 0.1872658121135139
$ws.Range("AE3").Value = 0.1872658121135139
$ws.Range("AF2").Value = 0.4732658121135139
$ws.Range("AF3").Value = 0.4732658121135139
$ws.Range("AG2").Value = -0.2757341878864861
$ws.Range("AG3").Value = -0.2757341878864861
$ws.Range("AH2").Value = 0.1700397461333858
$ws.Range("AH3").Value = 0.1700397461333858
$ws.Range("AI2").Value = 0.3433777490190954
$ws.Range("AI3").Value = 0.3433777490190954
$ws.Range("AJ2").Value = -0.135544817321592
$ws.Range("AJ3").Value = -0.135544817321592
$ws.Range("AK2").Value = -0.4381839638806663
$ws.Range("AK3").Value = -0.4381839638806663
$ws.Range("AL2").Value = 0.016
$ws.Range("AL3").Value = 0.016
$ws.Range("AM2").Value = 0.016
$ws.Range("AM3").Value = 0.016
$ws.Range("AN2").Value = -0.3466666666666667
$ws.Range("AN3").Value = -0.3466666666666667
$ws.Range("AO2").Value = -56.375
$ws.Range("AO3").Value = -56.375
$ws.Range("AP2").Value = 0.3342232580442256
$ws.Range("AP3").Value = 0.3342232580442256
$ws.Range("AQ2").Value = -56.375
$ws.Range("AQ3").Value = -56.375
